$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("meat")

# Add the new "food supergroup" property row to the meat sheet
$ws.Range("A10").Value = "food supergroup"
$ws.Range("B10").Value = "meat"

# Make the "meat" sheet the active tab and select B1 (matches the new
# selection/tabSelected state recorded in the saved workbook)
$ws.Activate()
$ws.Range("B1").Select()
